$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.712.84"
$ws.Range("E2").Value = "  +8.29%  "
$ws.Range("D3").Value = "1.774.44"
$ws.Range("E3").Value = "  +4.03%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'225.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'30.48"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.46%  "
$ws.Range("D9").Value = "'46.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("E10").Value = "  +3.59%  "
$ws.Range("D11").Value = "'0.0662"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.26%  "
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "2.029.26"
$ws.Range("E13").Value = "  +4.01%  "
$ws.Range("D14").Value = "1.771.75"
$ws.Range("E14").Value = "  +3.96%  "
$ws.Range("D15").Value = "'0.626"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "33.726.46"
$ws.Range("E16").Value = "  +8.09%  "
$ws.Range("D17").Value = "'9.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.53%  "
$ws.Range("D18").Value = "'4.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "'68.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D20").Value = "'251.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "'10.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").Value = "'4.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").Value = "'159.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "'16.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("D29").Value = "'6.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.64%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'3.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.48%  "
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("D34").Value = "'3.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.10%  "
$ws.Range("E35").Value = "  +5.09%  "
$ws.Range("D36").Value = "1.482.20"
$ws.Range("E36").Value = "  -2.44%  "
$ws.Range("E37").Value = "  +3.37%  "
$ws.Range("D38").Value = "'0.635"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.36%  "
$ws.Range("E39").Value = "  +2.47%  "
$ws.Range("D40").Value = "'83.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").Value = "'2.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.07%  "
$ws.Range("D42").Value = "'2.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").Value = "'0.884"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.96%  "
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("D45").Value = "'0.0512"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("D47").Value = "1.927.57"
$ws.Range("E47").Value = "  +4.77%  "
$ws.Range("D48").Value = "'5.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "'11.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.21%  "
$ws.Range("D51").Value = "'50.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.92%  "